# OM-distances2.xlsx — add a date series (E30:E47) used as the chart's
# category axis labels, name the series "Evolutionary rate", and turn on
# the (auto) chart title.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Write the new date labels into E30:E47 -----------------------------
$dates = @(
    "2011-06", "2011-07", "2011-08", "2011-09", "2011-10", "2011-11",
    "2011-12", "2012-01", "2012-02", "2012-03", "2012-04", "2012-05",
    "2012-06", "2012-07", "2012-08", "2012-09", "2012-10", "2012-11"
)
for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = 30 + $i
    $ws.Range("E$row").Value = $dates[$i]
}

# --- 2. Update the chart: series title + categories -------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$chart.HasTitle = $true

$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(""Evolutionary rate"",'OM-distances2.csv'!`$E`$30:`$E`$47,'OM-distances2.csv'!`$C`$25:`$T`$25,1)"

# --- 3. Restore the selection to match the edited range ---------------------
$ws.Range("E30:E47").Select()
